$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 1.044828013438889
$ws.Range("R2").Value = 9.40345212095
$ws.Range("S2").Value = 0.2087950866344732
$ws.Range("T2").Value = 0.2087950866344732

# Row 3
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("S3").Value = 0.4253229592313036
$ws.Range("T3").Value = 0.4253229592313036

# Row 4
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 1.473789855461111
$ws.Range("R4").Value = 13.26410869915
$ws.Range("S4").Value = 0.2945174484164121
$ws.Range("T4").Value = 0.2945174484164122

# Row 5
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 0.3571139337666666
$ws.Range("R5").Value = 3.2140254039
$ws.Range("S5").Value = 0.07136450571781097
$ws.Range("T5").Value = 0.07136450571781099

$wb.Save()
